$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data referencing the named ranges, with blank columns in between (C6, E6)
$ws.Range("A6").Value = "RangeOverBlanks"
$ws.Range("B6").Formula = "=B1"
$ws.Range("D6").Formula = "=B2"
$ws.Range("F6").Formula = "=B3"

# Define the new named range that spans over blank cells (B6:F6)
$wb.Names.Add("RangeOverBlanks", "=Sheet1!`$B`$6:`$F`$6")

# Update the selection to match the target state
$ws.Range("F7").Select()
